# The source document carries a stray `_GoBack` bookmark (Word's "last
# edit location" marker) sitting right after "Complaint details", before
# the trailing colon run. The target edit moves that marker so it instead
# sits inside the "Date:" field label, splitting "Date:" into "D" + "ate:"
# runs around it - i.e. the last edit was made there instead.

$d = $word.ActiveDocument

# --- 1. Drop the old `_GoBack` bookmark after "Complaint details". -------
# Bookmark.Delete() is a no-op against the existing (pre-loaded) bookmark
# in this runtime, but rewriting the text of a Range that spans the
# bookmark's position (via Find & Replace) does drop it - the paragraph's
# "Complaint details" + ":" runs get rebuilt without the bookmark between
# them.
$rng1 = $d.Content
[void]$rng1.Find.Execute("Complaint details:", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Complaint details:", 2)

# --- 2. Re-create `_GoBack` inside the "Date:" label. ---------------------
# Find the "Date:" run and add a fresh bookmark one character in (right
# after the "D"), which splits the run into "D" and "ate:" exactly as a
# real caret placement would.
$rng2 = $d.Content
[void]$rng2.Find.Execute("Date:", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$caret = $d.Range($rng2.Start + 1, $rng2.Start + 1)
$d.Bookmarks.Add("_GoBack", $caret)
